$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.101.46'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.468.62'
$ws.Range('E3').Value = '  -1.56%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.66'
$ws.Range('E5').Value = '  -1.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.00'
$ws.Range('E6').Value = '  -2.54%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.514'
$ws.Range('E8').Value = '  -1.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.467.67'
$ws.Range('E9').Value = '  -1.61%  '
$ws.Range('E10').Value = '  -2.34%  '
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.98'
$ws.Range('E12').Value = '  -2.56%  '
$ws.Range('E13').Value = '  -2.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.47'
$ws.Range('E14').Value = '  -3.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.916.69'
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.023.28'
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('E17').Value = '  -4.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.500.04'
$ws.Range('E18').Value = '  -1.86%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.15'
$ws.Range('E19').Value = '  -5.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.54'
$ws.Range('E20').Value = '  -3.91%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '352.93'
$ws.Range('E21').Value = '  -3.76%  '
$ws.Range('E22').Value = '  -2.74%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.15'
$ws.Range('E24').Value = '  -3.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.22'
$ws.Range('E25').Value = '  -7.58%  '
$ws.Range('E26').Value = '  -7.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.20'
$ws.Range('E27').Value = '  -7.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.593.86'
$ws.Range('E29').Value = '  -0.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0902'
$ws.Range('E30').Value = '  -6.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '511.66'
$ws.Range('E31').Value = '  -4.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.72'
$ws.Range('E32').Value = '  -7.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.78'
$ws.Range('E33').Value = '  -4.98%  '
$ws.Range('E34').Value = '  -6.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.120'
$ws.Range('E36').Value = '  -6.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '159.67'
$ws.Range('E37').Value = '  +0.78%  '
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('E39').Value = '  -2.51%  '
$ws.Range('E40').Value = '  -5.39%  '
$ws.Range('E41').Value = '  -0.42%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.66'
$ws.Range('E42').Value = '  -6.32%  '
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.326'
$ws.Range('E43').Value = '  -6.57%  '
$ws.Range('E44').Value = '  -6.43%  '
$ws.Range('E45').Value = '  -4.75%  '
$ws.Range('E46').Value = '  -2.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '140.76'
$ws.Range('E47').Value = '  -3.57%  '
$ws.Range('E48').Value = '  -6.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.514'
$ws.Range('E49').Value = '  -6.36%  '
$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.59'
$ws.Range('E50').Value = '  -6.89%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0252'
$ws.Range('E51').Value = '  -10.13%  '
